$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the cells that were removed from the test data (C3:D5)
$ws.Range("C3:D5").ClearContents()

# Update the current selection to match the new state of the sheet
[void]$ws.Range("C3:D5").Select()
